$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SQL Parser")

# --- Step 1: copy formatting (fills) between cells BEFORE changing values ---
# C7 picks up the "expressions-style" fill (same as C2 / style index 9)
$ws.Range("C2").Copy()
$ws.Range("C7").PasteSpecial(-4122)   # xlPasteFormats

# F2 picks up the fill currently on F5 (style index 9)
$ws.Range("F5").Copy()
$ws.Range("F2").PasteSpecial(-4122)   # xlPasteFormats

# F5 loses its fill entirely (back to default/no style)
$ws.Range("F5").ClearFormats()

$excel.CutCopyMode = 0

# --- Step 2: update cell text content (column F, rows 2-5) ---
$ws.Range("F2").Value = "print RA tree"
$ws.Range("F3").Value = "like"
$ws.Range("F4").Value = "views"
$ws.Range("F5").Value = "substring"

# --- Step 3: update cell text content (column C, rows 7-10) ---
$ws.Range("C7").Value = "parenthesis"
$ws.Range("C8").Value = "exists"
$ws.Range("C9").Value = "in"
$ws.Range("C10").Value = "between"

# --- Step 4: update the active selection on the sheet ---
$ws.Range("G4").Select()
